# Update the LR-pair TPM-derived metrics for Ccl11-Ccr3 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs)
$ws.Range("G2").Value = 1.57077
$ws.Range("H2").Value = 4.71231
$ws.Range("I2").Value = 0.02582502173444737
$ws.Range("J2").Value = 0.02582502173444737
$ws.Range("M2").Value = 0.08962966666666666
$ws.Range("N2").Value = 0.268889
$ws.Range("O2").Value = 0.4339761198462219
$ws.Range("P2").Value = 0.4339761198462219
$ws.Range("Q2").Value = 0.14078759151
$ws.Range("R2").Value = 1.26708832359
$ws.Range("S2").Value = 0.01120744272725982
$ws.Range("T2").Value = 0.01120744272725982

# Row 3 (ECs -> MuSCs)
$ws.Range("G3").Value = 1.57077
$ws.Range("H3").Value = 4.71231
$ws.Range("I3").Value = 0.02582502173444737
$ws.Range("J3").Value = 0.02582502173444737
$ws.Range("M3").Value = 0.1169016666666667
$ws.Range("N3").Value = 0.350705
$ws.Range("O3").Value = 0.5660238801537781
$ws.Range("P3").Value = 0.5660238801537781
$ws.Range("Q3").Value = 0.18362563095
$ws.Range("R3").Value = 1.65263067855
$ws.Range("S3").Value = 0.01461757900718755
$ws.Range("T3").Value = 0.01461757900718755

# Row 4 (FAPs -> FAPs)
$ws.Range("I4").Value = 0.934831682683009
$ws.Range("J4").Value = 0.934831682683009
$ws.Range("M4").Value = 0.08962966666666666
$ws.Range("N4").Value = 0.268889
$ws.Range("O4").Value = 0.4339761198462219
$ws.Range("P4").Value = 0.4339761198462219
$ws.Range("Q4").Value = 5.096324890856779
$ws.Range("R4").Value = 45.86692401771101
$ws.Range("S4").Value = 0.4056946263600867
$ws.Range("T4").Value = 0.4056946263600867

# Row 5 (FAPs -> MuSCs)
$ws.Range("I5").Value = 0.934831682683009
$ws.Range("J5").Value = 0.934831682683009
$ws.Range("M5").Value = 0.1169016666666667
$ws.Range("N5").Value = 0.350705
$ws.Range("O5").Value = 0.5660238801537781
$ws.Range("P5").Value = 0.5660238801537781
$ws.Range("Q5").Value = 6.647005347366112
$ws.Range("R5").Value = 59.82304812629501
$ws.Range("S5").Value = 0.5291370563229222
$ws.Range("T5").Value = 0.5291370563229222

# Row 6 (Inflammatory-Mac -> FAPs)
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.9273763333333335
$ws.Range("H6").Value = 2.782129
$ws.Range("I6").Value = 0.01524698967025436
$ws.Range("J6").Value = 0.01524698967025436
$ws.Range("M6").Value = 0.08962966666666666
$ws.Range("N6").Value = 0.268889
$ws.Range("O6").Value = 0.4339761198462219
$ws.Range("P6").Value = 0.4339761198462219
$ws.Range("Q6").Value = 0.08312043163122224
$ws.Range("R6").Value = 0.7480838846810001
$ws.Range("S6").Value = 0.006616829416432413
$ws.Range("T6").Value = 0.006616829416432413

# Row 7 (Inflammatory-Mac -> MuSCs)
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.9273763333333335
$ws.Range("H7").Value = 2.782129
$ws.Range("I7").Value = 0.01524698967025436
$ws.Range("J7").Value = 0.01524698967025436
$ws.Range("M7").Value = 0.1169016666666667
$ws.Range("N7").Value = 0.350705
$ws.Range("O7").Value = 0.5660238801537781
$ws.Range("P7").Value = 0.5660238801537781
$ws.Range("Q7").Value = 0.1084118389938889
$ws.Range("R7").Value = 0.975706550945
$ws.Range("S7").Value = 0.008630160253821947
$ws.Range("T7").Value = 0.008630160253821947

# Row 8 (MuSCs -> FAPs)
$ws.Range("G8").Value = 0.7810079999999999
$ws.Range("H8").Value = 2.343024
$ws.Range("I8").Value = 0.0128405486320577
$ws.Range("J8").Value = 0.0128405486320577
$ws.Range("M8").Value = 0.08962966666666666
$ws.Range("N8").Value = 0.268889
$ws.Range("O8").Value = 0.4339761198462219
$ws.Range("P8").Value = 0.4339761198462219
$ws.Range("Q8").Value = 0.07000148670399998
$ws.Range("R8").Value = 0.630013380336
$ws.Range("S8").Value = 0.005572491472037111
$ws.Range("T8").Value = 0.005572491472037111

# Row 9 (MuSCs -> MuSCs)
$ws.Range("G9").Value = 0.7810079999999999
$ws.Range("H9").Value = 2.343024
$ws.Range("I9").Value = 0.0128405486320577
$ws.Range("J9").Value = 0.0128405486320577
$ws.Range("M9").Value = 0.1169016666666667
$ws.Range("N9").Value = 0.350705
$ws.Range("O9").Value = 0.5660238801537781
$ws.Range("P9").Value = 0.5660238801537781
$ws.Range("Q9").Value = 0.09130113687999999
$ws.Range("R9").Value = 0.8217102319199999
$ws.Range("S9").Value = 0.007268057160020585
$ws.Range("T9").Value = 0.007268057160020585

# Row 10 (Resolving-Mac -> FAPs)
$ws.Range("G10").Value = 0.6846153333333334
$ws.Range("H10").Value = 2.053846
$ws.Range("I10").Value = 0.01125575728023152
$ws.Range("J10").Value = 0.01125575728023152
$ws.Range("M10").Value = 0.08962966666666666
$ws.Range("N10").Value = 0.268889
$ws.Range("O10").Value = 0.4339761198462219
$ws.Range("P10").Value = 0.4339761198462219
$ws.Range("Q10").Value = 0.06136184412155556
$ws.Range("R10").Value = 0.552256597094
$ws.Range("S10").Value = 0.004884729870405738
$ws.Range("T10").Value = 0.004884729870405738

# Row 11 (Resolving-Mac -> MuSCs)
$ws.Range("G11").Value = 0.6846153333333334
$ws.Range("H11").Value = 2.053846
$ws.Range("I11").Value = 0.01125575728023152
$ws.Range("J11").Value = 0.01125575728023152
$ws.Range("M11").Value = 0.1169016666666667
$ws.Range("N11").Value = 0.350705
$ws.Range("O11").Value = 0.5660238801537781
$ws.Range("P11").Value = 0.5660238801537781
$ws.Range("Q11").Value = 0.08003267349222222
$ws.Range("R11").Value = 0.72029406143
$ws.Range("S11").Value = 0.006371027409825781
$ws.Range("T11").Value = 0.006371027409825781
